$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Greek-letter notations introduced for the organism-level community model
$ws.Range("G19").Value = "η"
$ws.Range("G22").Value = "υ"

# G20 re-uses the existing "λ" notation (ψ is retired); G21/G23 keep their existing values
$ws.Range("G20").Value = "λ"
$ws.Range("G21").Value = "κ"
$ws.Range("G23").Value = "ω"

# New Meaning column (H) describing each parameter
$ws.Range("H21").Value = "organism-level community identity probability"
$ws.Range("H19").Value = "concentration parameter describing the C-dimensional Dirichlet distribution of organism-level community probabilities within a sample "
$ws.Range("H20").Value = "probability describing the S-dimensional categorical distribution of species identity of an organism, given its community identity"
$ws.Range("H22").Value = "concentration parameter describing organism-level community identity probability (the community proportion)"
$ws.Range("H23").Value = "probability defining the categorical distribution of community identity of individual organisms"

# Update the selected cell to reflect the new last-edited cell (H23)
$ws.Range("H23").Select()
